# Auto-generated cell updates for cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers need to be forced to
# Text format before assignment so Excel doesn't silently convert them into
# numeric values (which would lose formatting like trailing zeros).
# The style is reset back to Normal afterwards so no stray formatting is left
# behind on the cell (only the string-vs-number storage type matters).

$ws.Range("D2").Value = '30.932.31'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '1.908.37'
$ws.Range("E3").Value = '  +0.73%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4905'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2967'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06774'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").Value = '1.907.33'
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.10'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07285'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '89.96'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.119'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6712'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.90%  '
$ws.Range("D16").Value = '30.897.44'
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007952'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.50'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9992'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = '2.152.24'
$ws.Range("E20").Value = '  +0.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.093'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '207.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.210'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.649'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.968'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.429'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.318'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09176'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.043'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05174'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7510'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.116'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.705'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01839'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.726'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.113'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9270'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4486'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '106.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.824'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("E44").Value = '  +0.67%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.769'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1374'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '66.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +14.34%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.983'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.84%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.79%  '
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4072'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.44%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05912'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.91%  '
